# Updated recomputed TPM-derived statistics (ligand/receptor expression,
# specificity, and edge-weight columns) for rows 2-10 of the LR-pairs sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.433768000000001
$ws.Range("H2").Value = 16.301304
$ws.Range("I2").Value = 0.1262505823713576
$ws.Range("J2").Value = 0.1262505823713576
$ws.Range("M2").Value = 10.61298733333333
$ws.Range("N2").Value = 31.838962
$ws.Range("O2").Value = 0.1371494889257481
$ws.Range("P2").Value = 0.1371494889257481
$ws.Range("Q2").Value = 57.668510956272
$ws.Range("R2").Value = 519.016598606448
$ws.Range("S2").Value = 0.01731520284880977
$ws.Range("T2").Value = 0.01731520284880976

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.433768000000001
$ws.Range("H3").Value = 16.301304
$ws.Range("I3").Value = 0.1262505823713576
$ws.Range("J3").Value = 0.1262505823713576
$ws.Range("M3").Value = 42.26455300000001
$ws.Range("O3").Value = 0.5461762707865787
$ws.Range("P3").Value = 0.5461762707865786
$ws.Range("Q3").Value = 229.6557756257041
$ws.Range("R3").Value = 2066.901980631337
$ws.Range("S3").Value = 0.06895507226422189
$ws.Range("T3").Value = 0.06895507226422186

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.433768000000001
$ws.Range("H4").Value = 16.301304
$ws.Range("I4").Value = 0.1262505823713576
$ws.Range("J4").Value = 0.1262505823713576
$ws.Range("M4").Value = 24.50508366666667
$ws.Range("N4").Value = 73.51525100000001
$ws.Range("O4").Value = 0.3166742402876732
$ws.Range("P4").Value = 0.3166742402876732
$ws.Range("Q4").Value = 133.154939465256
$ws.Range("R4").Value = 1198.394455187304
$ws.Range("S4").Value = 0.03998030725832599
$ws.Range("T4").Value = 0.03998030725832597

# Row 5
$ws.Range("G5").Value = 5.270503666666666
$ws.Range("I5").Value = 0.1224572262391479
$ws.Range("J5").Value = 0.1224572262391479
$ws.Range("M5").Value = 10.61298733333333
$ws.Range("N5").Value = 31.838962
$ws.Range("O5").Value = 0.1371494889257481
$ws.Range("P5").Value = 0.1371494889257481
$ws.Range("Q5").Value = 55.93578865462022
$ws.Range("R5").Value = 503.422097891582
$ws.Range("S5").Value = 0.01679494599396385
$ws.Range("T5").Value = 0.01679494599396385

# Row 6
$ws.Range("G6").Value = 5.270503666666666
$ws.Range("I6").Value = 0.1224572262391479
$ws.Range("J6").Value = 0.1224572262391479
$ws.Range("M6").Value = 42.26455300000001
$ws.Range("O6").Value = 0.5461762707865787
$ws.Range("P6").Value = 0.5461762707865786
$ws.Range("S6").Value = 0.06688323115816619
$ws.Range("T6").Value = 0.06688323115816616

# Row 7
$ws.Range("G7").Value = 5.270503666666666
$ws.Range("I7").Value = 0.1224572262391479
$ws.Range("J7").Value = 0.1224572262391479
$ws.Range("M7").Value = 24.50508366666667
$ws.Range("N7").Value = 73.51525100000001
$ws.Range("O7").Value = 0.3166742402876732
$ws.Range("P7").Value = 0.3166742402876732
$ws.Range("Q7").Value = 129.1541333171401
$ws.Range("R7").Value = 1162.387199854261
$ws.Range("S7").Value = 0.03877904908701788
$ws.Range("T7").Value = 0.03877904908701788

# Row 8
$ws.Range("G8").Value = 32.33527633333333
$ws.Range("H8").Value = 97.00582900000001
$ws.Range("I8").Value = 0.7512921913894945
$ws.Range("J8").Value = 0.7512921913894944
$ws.Range("M8").Value = 10.61298733333333
$ws.Range("N8").Value = 31.838962
$ws.Range("O8").Value = 0.1371494889257481
$ws.Range("P8").Value = 0.1371494889257481
$ws.Range("Q8").Value = 343.1738781454998
$ws.Range("R8").Value = 3088.564903309498
$ws.Range("S8").Value = 0.1030393400829745
$ws.Range("T8").Value = 0.1030393400829745

# Row 9
$ws.Range("G9").Value = 32.33527633333333
$ws.Range("H9").Value = 97.00582900000001
$ws.Range("I9").Value = 0.7512921913894945
$ws.Range("J9").Value = 0.7512921913894944
$ws.Range("M9").Value = 42.26455300000001
$ws.Range("O9").Value = 0.5461762707865787
$ws.Range("P9").Value = 0.5461762707865786
$ws.Range("Q9").Value = 1366.636000359812
$ws.Range("R9").Value = 12299.72400323831
$ws.Range("S9").Value = 0.4103379673641907
$ws.Range("T9").Value = 0.4103379673641905

# Row 10
$ws.Range("G10").Value = 32.33527633333333
$ws.Range("H10").Value = 97.00582900000001
$ws.Range("I10").Value = 0.7512921913894945
$ws.Range("J10").Value = 0.7512921913894944
$ws.Range("M10").Value = 24.50508366666667
$ws.Range("N10").Value = 73.51525100000001
$ws.Range("O10").Value = 0.3166742402876732
$ws.Range("P10").Value = 0.3166742402876732
$ws.Range("Q10").Value = 792.37865193312
$ws.Range("R10").Value = 7131.40786739808
$ws.Range("S10").Value = 0.2379148839423293
$ws.Range("T10").Value = 0.2379148839423293
